{"js": "// The paragraph mentioning the author's name/NIM inside the printed Java\n// source code is updated to the actual author (Ajeng Nur Khorisa /\n// E31201199), and the stray \"_GoBack\" bookmark that Word drops at the\n// last edited location is moved from the title paragraph down to right\n// after the new NIM text.\n\nconst body = context.document.body;\n\n// 1) Drop the \"_GoBack\" bookmark from its old spot (end of the title\n//    paragraph \"Laporan Praktikum ke 2\").\ncontext.document.deleteBookmark(\"_GoBack\");\n\n// 2) Replace the placeholder name \"Lukman Afandi\" with \"Ajeng Nur Khorisa\"\n//    inside `System.out.println(\"Nama: Lukman Afandi\");`.\nconst nameResults = body.search(\"Lukman Afandi\", { matchCase: true });\nnameResults.load(\"text\");\nawait context.sync();\n\nif (nameResults.items.length > 0) {\n  nameResults.items[0].insertText(\"Ajeng Nur Khorisa\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) Replace the placeholder NIM \"E31200844\" with \"E31201199\" inside\n//    `System.out.println(\"NIM: E31200844\");`.\nconst nimResults = body.search(\"E31200844\", { matchCase: true });\nnimResults.load(\"text\");\nawait context.sync();\n\nif (nimResults.items.length > 0) {\n  const nimRange = nimResults.items[0];\n  nimRange.insertText(\"E31201199\", Word.InsertLocation.replace);\n  await context.sync();\n\n  // 4) Re-plant \"_GoBack\" immediately after the freshly typed NIM value \u2014\n  //    the same spot Word leaves it after the last text edit.\n  const endOfNim = nimRange.getRange(Word.RangeLocation.end);\n  endOfNim.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# The paragraph mentioning the author's name/NIM inside the printed Java\n# source code is updated to the actual author (Ajeng Nur Khorisa /\n# E31201199), and the stray \"_GoBack\" bookmark that Word drops at the\n# last edited location is moved from the title paragraph down to right\n# after the new NIM text.\n\n$d = $word.ActiveDocument\n\n# 1) Drop the \"_GoBack\" bookmark from its old spot (end of the title\n#    paragraph \"Laporan Praktikum ke 2\").\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Replace the placeholder name \"Lukman Afandi\" with \"Ajeng Nur Khorisa\"\n#    inside `System.out.println(\"Nama: Lukman Afandi\");`.\n$wdFindContinue = 1\n$wdReplaceAll = 2\n$nameRange = $d.Content\n$nameRange.Find.Execute(\"Lukman Afandi\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"Ajeng Nur Khorisa\", $wdReplaceAll)\n\n# 3) Replace the placeholder NIM \"E31200844\" with \"E31201199\" inside\n#    `System.out.println(\"NIM: E31200844\");`.\n$nimRange = $d.Content\n$nimRange.Find.Execute(\"E31200844\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"E31201199\", $wdReplaceAll)\n\n# 4) Re-plant \"_GoBack\" immediately after the freshly typed NIM value \u2014\n#    the same spot Word leaves it after the last text edit.\n$locateRange = $d.Content\n$locateRange.Find.Execute(\"NIM: E31201199\", $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"\", 0)\n$bookmarkRange = $d.Range($locateRange.End, $locateRange.End)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n"}
